$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 16.75560066666667
$ws.Range("H2").Value = 50.266802
$ws.Range("I2").Value = 0.9064438825950116
$ws.Range("J2").Value = 0.9064438825950115
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 29.546424
$ws.Range("N2").Value = 88.63927200000001
$ws.Range("O2").Value = 0.9033225104610835
$ws.Range("P2").Value = 0.9033225104610834
$ws.Range("Q2").Value = 495.068081672016
$ws.Range("R2").Value = 4455.612735048144
$ws.Range("S2").Value = 0.8188111636178175
$ws.Range("T2").Value = 0.8188111636178174
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 16.75560066666667
$ws.Range("H3").Value = 50.266802
$ws.Range("I3").Value = 0.9064438825950116
$ws.Range("J3").Value = 0.9064438825950115
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.672785333333334
$ws.Range("N3").Value = 5.018356000000001
$ws.Range("O3").Value = 0.05114204841740398
$ws.Range("P3").Value = 0.05114204841740398
$ws.Range("Q3").Value = 28.02852304639023
$ws.Range("R3").Value = 252.256707417512
$ws.Range("S3").Value = 0.04635739693133373
$ws.Range("T3").Value = 0.04635739693133373
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 16.75560066666667
$ws.Range("H4").Value = 50.266802
$ws.Range("I4").Value = 0.9064438825950116
$ws.Range("J4").Value = 0.9064438825950115
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 1.489401
$ws.Range("N4").Value = 4.468203
$ws.Range("O4").Value = 0.04553544112151264
$ws.Range("P4").Value = 0.04553544112151264
$ws.Range("Q4").Value = 24.955808388534
$ws.Range("R4").Value = 224.602275496806
$ws.Range("S4").Value = 0.04127532204586047
$ws.Range("T4").Value = 0.04127532204586046
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 0.07688299999999999
$ws.Range("H5").Value = 0.230649
$ws.Range("I5").Value = 0.004159213770485276
$ws.Range("J5").Value = 0.004159213770485276
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 29.546424
$ws.Range("N5").Value = 88.63927200000001
$ws.Range("O5").Value = 0.9033225104610835
$ws.Range("P5").Value = 0.9033225104610834
$ws.Range("Q5").Value = 2.271617716392
$ws.Range("R5").Value = 20.444559447528
$ws.Range("S5").Value = 0.003757111424699068
$ws.Range("T5").Value = 0.003757111424699068
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 0.07688299999999999
$ws.Range("H6").Value = 0.230649
$ws.Range("I6").Value = 0.004159213770485276
$ws.Range("J6").Value = 0.004159213770485276
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.672785333333334
$ws.Range("N6").Value = 5.018356000000001
$ws.Range("O6").Value = 0.05114204841740398
$ws.Range("P6").Value = 0.05114204841740398
$ws.Range("Q6").Value = 0.1286087547826667
$ws.Range("R6").Value = 1.157478793044
$ws.Range("S6").Value = 0.0002127107120284914
$ws.Range("T6").Value = 0.0002127107120284914
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 0.07688299999999999
$ws.Range("H7").Value = 0.230649
$ws.Range("I7").Value = 0.004159213770485276
$ws.Range("J7").Value = 0.004159213770485276
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.489401
$ws.Range("N7").Value = 4.468203
$ws.Range("O7").Value = 0.04553544112151264
$ws.Range("P7").Value = 0.04553544112151264
$ws.Range("Q7").Value = 0.114509617083
$ws.Range("R7").Value = 1.030586553747
$ws.Range("S7").Value = 0.0001893916337577169
$ws.Range("T7").Value = 0.0001893916337577169
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.652500333333333
$ws.Range("H8").Value = 4.957501
$ws.Range("I8").Value = 0.0893969036345032
$ws.Range("J8").Value = 0.08939690363450319
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 29.546424
$ws.Range("N8").Value = 88.63927200000001
$ws.Range("O8").Value = 0.9033225104610835
$ws.Range("P8").Value = 0.9033225104610834
$ws.Range("Q8").Value = 48.825475508808
$ws.Range("R8").Value = 439.429279579272
$ws.Range("S8").Value = 0.08075423541856699
$ws.Range("T8").Value = 0.08075423541856697
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.652500333333333
$ws.Range("H9").Value = 4.957501
$ws.Range("I9").Value = 0.0893969036345032
$ws.Range("J9").Value = 0.08939690363450319
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 1.672785333333334
$ws.Range("N9").Value = 5.018356000000001
$ws.Range("O9").Value = 0.05114204841740398
$ws.Range("P9").Value = 0.05114204841740398
$ws.Range("Q9").Value = 2.764278320928445
$ws.Range("R9").Value = 24.878504888356
$ws.Range("S9").Value = 0.00457194077404176
$ws.Range("T9").Value = 0.00457194077404176
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 1.652500333333333
$ws.Range("H10").Value = 4.957501
$ws.Range("I10").Value = 0.0893969036345032
$ws.Range("J10").Value = 0.08939690363450319
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.489401
$ws.Range("N10").Value = 4.468203
$ws.Range("O10").Value = 0.04553544112151264
$ws.Range("P10").Value = 0.04553544112151264
$ws.Range("Q10").Value = 2.461235648967
$ws.Range("R10").Value = 2.461235648967
$ws.Range("S10").Value = 0.00407072744189446
$ws.Range("T10").Value = 0.004070727441894459

Write-Output "Applied NATMI Jag2-Notch4 updates"
